$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Product ID values for rows 8-13 (continuing the
# sequence 1..6 used in rows 2-7) -- these correspond to the 3 new
# products / 2 new services added to the pricing sheet.
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12

# Move the active selection to A2
$ws.Range("A2").Select()
